$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, copying the same formatting as the
# existing header cells (e.g. G1: bold, bordered, centered).
$g1 = $ws.Range("G1")
$h1 = $ws.Range("H1")
$g1.Copy()
$h1.PasteSpecial(-4122)  # xlPasteFormats
$h1.Value = "Save"

# Fill in the Save column values for rows 2-7
$values = @(0, 0, 0, 1, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
